$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet1: introduce "Exorcism" (and related) as a new Spell Type (col D),
# shifting Binding/Stealth/Spirit Control/Barrier from Element (col E) into
# Spell Type (col D), and backfilling Element with new "Bind"/"Drain" entries.

$ws1.Range("D8").Value = "Exorcism"
$ws1.Range("D9").Value = "Binding"
$ws1.Range("D10").Value = "Stealth"
$ws1.Range("D11").Value = "Spirit Control"
$ws1.Range("D12").Value = "Barrier"

$ws1.Range("E10").Value = "Bind"
$ws1.Range("E11").Value = "Lightning"
$ws1.Range("E12").Value = "Drain"
$ws1.Range("E13").ClearContents()
$ws1.Range("E14").ClearContents()

# --- Sheet2: try out the new "Exorcism" spell-type combo
$ws2.Range("B3").Value = "Hand Seal"
$ws2.Range("C3").Value = "Exorcism"
$ws2.Range("D3").Value = "Drain"
$ws2.Range("E3").Value = "Roll Add"
$ws2.Range("F3").Value = "None"

# Update the selections to reflect where the author was working
$ws1.Range("E13").Select() | Out-Null
$ws2.Range("E3").Select() | Out-Null
